$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Child")

$ws.Range("D2").Value = "-10,8"
$ws.Range("D3").Value = "-6,1"
$ws.Range("D4").Value = "-1,8"
$ws.Range("D5").Value = "-1,0"
$ws.Range("D6").Value = "-6,-3"
$ws.Range("D7").Value = "0,-1"
$ws.Range("D8").Value = "6,-9"
$ws.Range("D9").Value = "5,-4"
$ws.Range("D10").Value = "5,1"
$ws.Range("D11").Value = "-10,6"
$ws.Range("D12").Value = "6,4"
$ws.Range("D13").Value = "4,1"
$ws.Range("D14").Value = "2,1"
$ws.Range("D15").Value = "4,3"
$ws.Range("D16").Value = "-2,-5"
$ws.Range("D17").Value = "-9,-1"
$ws.Range("D18").Value = "-8,2"
$ws.Range("D19").Value = "-4,7"
$ws.Range("D20").Value = "1,-6"
$ws.Range("D21").Value = "8,0"
